# This sheet (OOMII_DB) stores every value as text (the source data was
# exported with every cell, including blanks, written out as an inline string).
# The edit:
#   1. Cleans out the placeholder empty-string cells on rows 2 and 3, leaving
#      only the columns that actually carry a value.
#   2. Appends row 4: a new record with the same value pattern as row 2 (the
#      sparse / cleaned layout).
#   3. Appends row 5: a second new record, this one cloned in full from row 2 -
#      i.e. it still carries the blank placeholder cells in every other column,
#      mirroring how rows 2/3 looked before step 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Every cell we touch below should stay plain text (several values look like
# numbers, e.g. "0.0625" or "025", and must not be auto-converted).
$ws.Range("A2:OR5").NumberFormat = "@"

# Columns (1 = A ... 408 = OR) that are blank placeholders on rows 2 and 3.
$blankCols = @(
    2, 3, 4, 5, 6, 7, 14, 15, 16, 19, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 38, 45, 46, 47, 48, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 73, 74, 75, 78, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 97, 104, 105, 106, 107, 114, 115, 116, 117, 118, 119, 121, 122, 123, 124, 125, 126, 127, 128, 129, 130, 131, 132, 133, 134, 135, 136, 137, 138, 139, 146, 147, 148, 149, 151, 158, 159, 160, 161, 162, 163, 164, 165, 166, 167, 168, 169, 170, 171, 172, 173, 174, 175, 176, 177, 178, 179, 180, 187, 188, 189, 190, 191, 192, 193, 194, 195, 196, 197, 198, 199, 200, 201, 202, 203, 204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215, 216, 217, 218, 219, 220, 221, 222, 223, 224, 225, 226, 227, 228, 229, 230, 231, 238, 239, 240, 243, 250, 251, 252, 253, 254, 255, 256, 257, 258, 259, 260, 262, 269, 270, 271, 272, 279, 280, 281, 282, 283, 284, 285, 286, 287, 288, 289, 290, 291, 292, 293, 294, 295, 296, 297, 298, 299, 300, 301, 302, 303, 304, 305, 306, 307, 308, 309, 310, 311, 312, 313, 314, 315, 316, 317, 318, 319, 320, 321, 322, 323, 324, 325, 326, 327, 328, 329, 330, 331, 332, 333, 334, 335, 336, 337, 338, 339, 340, 341, 342, 343, 344, 345, 346, 347, 348, 349, 350, 351, 352, 353, 354, 355, 356, 357, 358, 359, 360, 361, 362, 363, 364, 365, 366, 367, 368, 369, 370, 371, 372, 373, 374, 375, 376, 377, 378, 379, 380, 381, 382, 383, 384, 385, 386, 387, 388, 389, 390, 391, 392, 393, 394, 395, 396, 397, 398, 399, 400, 401, 402, 403, 404, 405, 406, 407, 408
)

foreach ($c in $blankCols) {
    $ws.Cells.Item(2, $c).ClearContents()
    $ws.Cells.Item(3, $c).ClearContents()
}

# Row 2s populated values - reused verbatim for the new row 4 and row 5 records.
$rowValues = @{
    1 = "type01"
    8 = "0.0625"
    9 = "0.25"
    10 = "1.0"
    11 = "4.0"
    12 = "16.0"
    13 = "25.0"
    17 = "3.22"
    18 = "3.03"
    20 = "3.00"
    21 = "2.97"
    22 = "3.01"
    23 = "2.93"
    24 = "2.31"
    25 = "2.75"
    37 = "94"
    39 = "93"
    40 = "92"
    41 = "93"
    42 = "91"
    43 = "97"
    44 = "35"
    49 = ".1"
    50 = ".2"
    51 = ".1"
    52 = "4"
    53 = ".7"
    54 = ".9"
    67 = "0.0625"
    68 = "0.25"
    69 = "1.0"
    70 = "4.0"
    71 = "15.0"
    72 = "25.0"
    76 = "2.51"
    77 = "2.20"
    79 = "2.23"
    80 = "2.25"
    81 = "2.27"
    82 = "2.20"
    83 = "2.10"
    84 = "1.99"
    96 = "91"
    98 = "39"
    99 = "90"
    100 = "91"
    101 = "39"
    102 = "34"
    103 = "79"
    108 = ".2"
    109 = ".2"
    110 = ".1"
    111 = "4"
    112 = ".3"
    113 = ".13"
    120 = "220"
    140 = "2.57"
    141 = "1.75"
    142 = "1.67"
    143 = "1.69"
    144 = "1.73"
    145 = "1.69"
    150 = "68"
    152 = "55"
    153 = "66"
    154 = "63"
    155 = "60"
    156 = "60"
    157 = "51"
    181 = ".4"
    182 = ".3"
    183 = ".0"
    184 = ".3"
    185 = ".11"
    186 = ".24"
    232 = "0.0625"
    233 = "0.25"
    234 = "1.0"
    235 = "4.0"
    236 = "16.0"
    237 = "25.0"
    241 = "5.45"
    242 = "9.33"
    244 = "7.32"
    245 = "0.47"
    246 = "7.99"
    247 = "7.64"
    248 = "3.12"
    249 = "7.31"
    261 = "153"
    263 = "143"
    264 = "155"
    265 = "147"
    266 = "140"
    267 = "149"
    268 = "134"
    273 = ".9"
    274 = ".1"
    275 = ".7"
    276 = ".11"
    277 = ".5"
    278 = ".15"
}

foreach ($c in $rowValues.Keys) {
    $ws.Cells.Item(4, $c).Value = $rowValues[$c]
    $ws.Cells.Item(5, $c).Value = $rowValues[$c]
}

# Row 5 also keeps a (present-but-empty) cell for every other column, just like
# row 2 and row 3 used to before the clean-up above. A literal "" value isnt
# distinguishable from "no cell at all" through COM, so write it through a
# formula that evaluates to an empty string - the stored/display value is still
# an empty string, but the cell itself is materialised.
foreach ($c in $blankCols) {
    $ws.Cells.Item(5, $c).Formula = "="""""
}

Write-Output "done"
